$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 4743
$ws.Range("D2").Value = 81.70999999999999

$ws.Range("C3").Value = 1028
$ws.Range("D3").Value = 17.71

$ws.Range("C4").Value = 34
$ws.Range("D4").Value = 0.59

$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0
